# "Generate Report for Handback" - record that the fb3c528d... handback
# file has come back for both the zh-cn and de-de locales: stamp the
# handback file name + datetime, link the "Latest Target File" cell to the
# handback markdown on GitHub, and note (in Error Detail) that the handed
# back version is not the latest available.

$wb = $excel.ActiveWorkbook

$handbackMdDisplay = "fb3c528d-2209-4e1e-8dfe-30beac899105.md"
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6e70b82cbec75fcd8fefc9950d95d01a4858dc6/e2e/fb3c528d-2209-4e1e-8dfe-30beac899105.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ae3fcbff18256525f44cbd79d8a9e60607a32f9/e2e/fb3c528d-2209-4e1e-8dfe-30beac899105.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6e70b82cbec75fcd8fefc9950d95d01a4858dc6/e2e/fb3c528d-2209-4e1e-8dfe-30beac899105.md."

function Update-LocaleSheet {
    param([string]$SheetName, [string]$HandbackXlf, [string]$HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Error Detail column (P) is the widest-changed column: 13.75 -> 40 chars.
    $ws.Columns.Item(16).ColumnWidth = 39.16666666666667

    # I8 = "Latest Target File" -> link to the handback .md on GitHub.
    $ws.Hyperlinks.Add($ws.Range("I8"), $handbackUrl, "", "", $handbackMdDisplay) | Out-Null

    # J8 = "Latest Handback File"
    $ws.Range("J8").Value = $HandbackXlf

    # K8 = "Latest Handback DateTime"
    $ws.Range("K8").Value = $HandbackDateTime

    # P8 = "Error Detail"
    $ws.Range("P8").Value = $errorDetail
}

Update-LocaleSheet "zh-cn" "fb3c528d-2209-4e1e-8dfe-30beac899105.6a6082f863f96f33b0759c6e58eaccc52bd83913.zh-cn.xlf" "2016-09-01 00:47:15"
Update-LocaleSheet "de-de" "fb3c528d-2209-4e1e-8dfe-30beac899105.6a6082f863f96f33b0759c6e58eaccc52bd83913.de-de.xlf" "2016-09-01 00:47:22"
